# Weekly driver report update for 2025-04-21
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Bad Drivers table ---
# Row 3: Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.10.1
$ws.Range("C3").Value = 358
$ws.Range("D3").Value = 94.90000000000001

# Row 4: Intel(R) Dual Band Wireless-AC 8265 - 20.70.16.4
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 109
$ws.Range("D4").Value = 96.90000000000001

# Row 5: Intel(R) Wi-Fi 6 AX201 160MHz - 23.90.0.2
$ws.Range("C5").Value = 2188

# Row 6: Totals
$ws.Range("B6").Value = 32
$ws.Range("C6").Value = 2655

# --- Good Drivers table ---
# Row 14: Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4
$ws.Range("B14").Value = 449371

# Row 15: Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9
$ws.Range("B15").Value = 77999

# Row 23: Intel(R) Dual Band Wireless-AC 8265 - 20.70.12.5
$ws.Range("B23").Value = 144782
